# Applies the edits described by the commit diff:
#  - workbook-level: iterative-calc "max change" setting (iterateDelta), rename sheet
#  - sheet-level: best-fit-like column widths on A, C, F, G, H; move selection to F23

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook: iterative calculation "Maximum Change" (-> calcPr iterateDelta) ---
$excel.Iteration = $true
$excel.MaxChange = 0.0001
$excel.MaxIterations = 100

# --- Workbook: rename the sheet "Result" -> "1-15 Операторы" ---
$ws.Name = "1-15 Операторы"

# --- Worksheet: best-fit column widths (engine quantizes ColumnWidth to 1/6 char units,
#     so the COM "ColumnWidth" input is pre-adjusted to land as close as possible on the
#     target stored width from the diff) ---
$ws.Columns("A").ColumnWidth = 26.0
$ws.Columns("C").ColumnWidth = 9.5
$ws.Columns("F").ColumnWidth = 15.666666666666666
$ws.Columns("G").ColumnWidth = 11.833333333333334
$ws.Columns("H").ColumnWidth = 6.166666666666667

# --- Worksheet: move selection from E6 to F23 ---
$ws.Range("F23").Select()
